$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet1 "Trends Status" - updated counts/percentages
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 0
$ws1.Cells.Item(2, 4).Value = 0
$ws1.Cells.Item(2, 5).Value = 0

$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 4
$ws1.Cells.Item(3, 4).Value = 0
$ws1.Cells.Item(3, 5).Value = 21.1

$ws1.Cells.Item(4, 2).Value = 3
$ws1.Cells.Item(4, 3).Value = 12
$ws1.Cells.Item(4, 4).Value = 37.5
$ws1.Cells.Item(4, 5).Value = 63.2

$ws1.Cells.Item(5, 2).Value = 4
$ws1.Cells.Item(5, 3).Value = 1
$ws1.Cells.Item(5, 4).Value = 50
$ws1.Cells.Item(5, 5).Value = 5.3

$ws1.Cells.Item(6, 2).Value = 1
$ws1.Cells.Item(6, 3).Value = 2
$ws1.Cells.Item(6, 4).Value = 12.5
$ws1.Cells.Item(6, 5).Value = 10.5

$ws1.Cells.Item(7, 2).Value = 38
$ws1.Cells.Item(7, 3).Value = 135

$ws1.Cells.Item(8, 2).Value = 380
$ws1.Cells.Item(8, 3).Value = 272

# ------------------------------------------------------------------
# Sheet3 "Priority Status" - updated species counts
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")

$ws3.Cells.Item(2, 2).Value = 103
$ws3.Cells.Item(3, 2).Value = 286
$ws3.Cells.Item(4, 2).Value = 554

# ------------------------------------------------------------------
# Sheet4 "Species qualification" - renamed row label + new counts
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Cells.Item(2, 1).Value = "SoIB Assessment"
$ws4.Cells.Item(2, 2).Value = 426

$ws4.Cells.Item(3, 2).Value = 46
$ws4.Cells.Item(3, 3).Value = 8

$ws4.Cells.Item(4, 2).Value = 154
$ws4.Cells.Item(4, 3).Value = 19

# ------------------------------------------------------------------
# Sheet5 "High Priority break-up" -> rename to "Interannual update - High Pri"
# and refresh its values (interannual-update numbers)
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")

$ws5.Cells.Item(2, 2).Value = 68
$ws5.Cells.Item(2, 3).Value = 66
$ws5.Cells.Item(2, 4).Value = 68
$ws5.Cells.Item(2, 5).Value = 77.3

$ws5.Cells.Item(3, 2).Value = 35
$ws5.Cells.Item(3, 3).Value = 34
$ws5.Cells.Item(3, 4).Value = 20
$ws5.Cells.Item(3, 5).Value = 22.7

$ws5.Name = "Interannual update - High Pri"

# ------------------------------------------------------------------
# New Sheet6 "Major update - High Priority " - carries forward the
# original ("before") High Priority break-up figures
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "Major update - High Priority "

$ws6.Cells.Item(1, 1).Value = "Break-up"
$ws6.Cells.Item(1, 2).Value = "High Species (no.)"
$ws6.Cells.Item(1, 3).Value = "High Species (perc.)"
$ws6.Cells.Item(1, 4).Value = "New High Species (no.)"
$ws6.Cells.Item(1, 5).Value = "New High Species (perc.)"
$ws6.Range("A1:E1").Font.Bold = $true
$ws6.Range("A1:E1").HorizontalAlignment = -4108

$ws6.Cells.Item(2, 1).Value = "Trend New"
$ws6.Cells.Item(2, 2).Value = 1
$ws6.Cells.Item(2, 3).Value = 5
$ws6.Cells.Item(2, 4).Value = 1
$ws6.Cells.Item(2, 5).Value = 5

$ws6.Cells.Item(3, 1).Value = "IUCN"
$ws6.Cells.Item(3, 2).Value = 19
$ws6.Cells.Item(3, 3).Value = 95
$ws6.Cells.Item(3, 4).Value = 19
$ws6.Cells.Item(3, 5).Value = 95
